$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 11.02.2022 15:30"

# D2: delta price, was text "+0.6", now numeric 0.6
$ws.Range("D2").Value = 0.6

# E2: old date, was text "2022-02-11 15:15:03", now numeric date serial
# formatted the same way as the other rows' Old Datum column (E3:E10)
$ws.Range("E2").Value = 44603.63545138889
$ws.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
